# Updates the Price (D) and Volume(1h) (E) columns on the cryptos sheet
# to match the latest scraped values, per the GitHub Actions data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "62.316.22"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  -1.27%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.001.99"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  -1.67%  "

$ws.Range("E4").Value = "  -0.09%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "588.61"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +0.40%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "145.99"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -3.46%  "

$ws.Range("E7").Value = "  -0.06%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.526"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -2.05%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.998.61"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -1.78%  "

$ws.Range("E10").Value = "  -4.17%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.78"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -0.76%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.465"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +3.90%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000229"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -2.37%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "34.47"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -4.80%  "

$ws.Range("E15").Value = "  +1.96%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.494.42"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -1.75%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "7.06"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -1.00%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "62.181.55"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -1.44%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "2.997.93"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -1.75%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "457.42"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -4.31%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.03"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -1.64%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.688"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -2.36%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.40"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -1.41%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "81.86"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -0.69%  "

$ws.Range("E25").Value = "  -9.02%  "

$ws.Range("E26").Value = "  -3.79%  "

$ws.Range("E27").Value = "  +0.00%  "

$ws.Range("E28").Value = "  -8.25%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.00"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -0.04%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.64"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -1.09%  "

$ws.Range("E31").Value = "  -5.58%  "

$ws.Range("E32").Value = "  -4.38%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "27.67"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +0.32%  "

$ws.Range("E34").Value = "  -1.61%  "

$ws.Range("E35").Value = "  -1.70%  "

$ws.Range("E36").Value = "  -3.35%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.74"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -2.53%  "

$ws.Range("E38").Value = "  -4.75%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "9.18"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -0.53%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "50.25"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -0.43%  "

$ws.Range("E41").Value = "  +7.17%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.88"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -11.19%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "393.47"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -9.31%  "

$ws.Range("E44").Value = "  -1.18%  "

$ws.Range("E45").Value = "  -6.93%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.732.56"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -3.40%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "37.35"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -2.62%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "129.38"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +0.25%  "

$ws.Range("E50").Value = "  -0.66%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.19"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -0.62%  "
